# Update the quantity/price figures in column B on Sheet1 to reflect
# the new supplier/server numbers ("moving to new server").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 25500   # was 30000
$ws.Range("B3").Value = 2500    # was 2100
$ws.Range("B4").Value = 10500   # was 12000
$ws.Range("B5").Value = 6000    # was 6150
